$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 520
$newTimestamp = "2023-01-07 12:56:47"

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 15).Value = $newTimestamp
}

$ws.Cells.Item(278, 13).Value = "Naturaplan Bio Artischocke 1 Stück - Online kein Bestand 1.95 Schweizer Franken"
$ws.Cells.Item(366, 13).Value = "Naturaplan Bio Shiitake-Pilze ca. 100g - Online kein Bestand 3.20 Schweizer Franken"
